# Update (Analyze PO & Forecast)
$wb = $excel.ActiveWorkbook

# --- Sheet: Forecast Comparison (MyForecast column updates) ---
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")

$wsForecast.Range("D2").Value = 60
$wsForecast.Range("D3").Value = 54
$wsForecast.Range("D4").Value = 51
$wsForecast.Range("D5").Value = 52
$wsForecast.Range("D6").Value = 52
$wsForecast.Range("D7").Value = 54
$wsForecast.Range("D8").Value = 52
$wsForecast.Range("D9").Value = 52
$wsForecast.Range("D10").Value = 50
$wsForecast.Range("D11").Value = 50
$wsForecast.Range("D12").Value = 44
$wsForecast.Range("D13").Value = 36
$wsForecast.Range("D17").Value = 29

# --- Sheet: Summary (derived stats, stored as text like the source) ---
$wsSummary = $wb.Worksheets.Item("Summary")

function Set-TextValue($range, [string]$text) {
    # Prefix with an apostrophe so the numeric-looking text is kept as text
    # (matches the original inline-string cell type), then reset the style
    # so no stray "quote prefix" number format is left behind on the cell.
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

Set-TextValue $wsSummary.Range("B9")  "804"
Set-TextValue $wsSummary.Range("B10") "427"
Set-TextValue $wsSummary.Range("B11") "217"
Set-TextValue $wsSummary.Range("B12") "60"
Set-TextValue $wsSummary.Range("B14") "29"
Set-TextValue $wsSummary.Range("B15") "2025-05-11"
